$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45189 -> 45190) for every data row, from row 2 through row 185.
$ws.Range("C2:C185").Value = 45190
